$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell values per the diff
$ws.Range("B1").Value = "Ria123456"
$ws.Range("B2").Value = "Ria12345"
$ws.Range("A3").Value = "ria12"
$ws.Range("B3").Value = "Ria123456"

# Remove row 5 entirely (its data is deleted from the sheet)
$ws.Range("A5:B5").EntireRow.Delete()

# Update selection to match target view state
$ws.Range("B5").Select()
